$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe preserves the existing "quote prefix" text-entry style
# (keeps cell format index s="4") while storing the text without the apostrophe.
$ws.Range("M2").Value = "'COMPLETE"
$ws.Range("M3").Value = "'COMPLETE"
